$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '66.094.60'
Set-TextValue $ws.Range("E2") '  +1.17%  '
Set-TextValue $ws.Range("D3") '2.697.49'
Set-TextValue $ws.Range("E3") '  +1.97%  '
Set-TextValue $ws.Range("E4") '  -0.03%  '
Set-TextValue $ws.Range("D5") '612.35'
Set-TextValue $ws.Range("E5") '  +1.43%  '
Set-TextValue $ws.Range("D6") '158.51'
Set-TextValue $ws.Range("E6") '  +1.51%  '
Set-TextValue $ws.Range("E7") '  -0.02%  '
Set-TextValue $ws.Range("D8") '0.591'
Set-TextValue $ws.Range("E8") '  +0.56%  '
Set-TextValue $ws.Range("E9") '  +5.50%  '
Set-TextValue $ws.Range("D10") '6.06'
Set-TextValue $ws.Range("E10") '  +4.62%  '
Set-TextValue $ws.Range("D11") '0.404'
Set-TextValue $ws.Range("E11") '  -0.75%  '
Set-TextValue $ws.Range("E12") '  +0.32%  '
Set-TextValue $ws.Range("E13") '  +11.10%  '
Set-TextValue $ws.Range("D14") '30.18'
Set-TextValue $ws.Range("E14") '  +2.98%  '
Set-TextValue $ws.Range("D15") '3.183.52'
Set-TextValue $ws.Range("E15") '  +2.04%  '
Set-TextValue $ws.Range("D16") '65.962.03'
Set-TextValue $ws.Range("E16") '  +1.19%  '
Set-TextValue $ws.Range("D17") '2.693.45'
Set-TextValue $ws.Range("E17") '  +2.78%  '
Set-TextValue $ws.Range("D18") '12.79'
Set-TextValue $ws.Range("E18") '  +0.42%  '
Set-TextValue $ws.Range("D19") '4.90'
Set-TextValue $ws.Range("E19") '  +0.05%  '
Set-TextValue $ws.Range("D20") '7.79'
Set-TextValue $ws.Range("E20") '  +5.70%  '
Set-TextValue $ws.Range("D21") '359.46'
Set-TextValue $ws.Range("E21") '  -0.03%  '
Set-TextValue $ws.Range("D22") '71.43'
Set-TextValue $ws.Range("E22") '  +3.38%  '
Set-TextValue $ws.Range("D23") '1.00'
Set-TextValue $ws.Range("E23") '  -0.08%  '
Set-TextValue $ws.Range("D24") '0.0000114'
Set-TextValue $ws.Range("E24") '  +18.27%  '
Set-TextValue $ws.Range("D25") '9.96'
Set-TextValue $ws.Range("E25") '  +5.40%  '
Set-TextValue $ws.Range("E26") '  -2.36%  '
Set-TextValue $ws.Range("D27") '1.68'
Set-TextValue $ws.Range("E27") '  +1.14%  '
Set-TextValue $ws.Range("E28") '  +3.82%  '
Set-TextValue $ws.Range("D29") '8.34'
Set-TextValue $ws.Range("E29") '  +1.17%  '
Set-TextValue $ws.Range("E30") '  +0.46%  '
Set-TextValue $ws.Range("E31") '  +0.05%  '
Set-TextValue $ws.Range("D32") '538.01'
Set-TextValue $ws.Range("E32") '  -0.63%  '
Set-TextValue $ws.Range("E33") '  -0.43%  '
Set-TextValue $ws.Range("D34") '6.70'
Set-TextValue $ws.Range("E34") '  +4.85%  '
Set-TextValue $ws.Range("D35") '5.51'
Set-TextValue $ws.Range("E35") '  -0.84%  '
Set-TextValue $ws.Range("E36") '  +1.54%  '
Set-TextValue $ws.Range("D37") '20.79'
Set-TextValue $ws.Range("E37") '  +0.50%  '
Set-TextValue $ws.Range("D38") '163.95'
Set-TextValue $ws.Range("E38") '  +1.31%  '
Set-TextValue $ws.Range("D39") '2.00'
Set-TextValue $ws.Range("E39") '  -0.87%  '
Set-TextValue $ws.Range("E40") '  -0.01%  '
Set-TextValue $ws.Range("D41") '0.999'
Set-TextValue $ws.Range("E41") '  +0.03%  '
Set-TextValue $ws.Range("B42") 'Aave'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D42") '168.64'
Set-TextValue $ws.Range("E42") '  +1.55%  '
Set-TextValue $ws.Range("B43") 'OKB'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D43") '42.63'
Set-TextValue $ws.Range("E43") '  -0.15%  '
Set-TextValue $ws.Range("D44") '4.18'
Set-TextValue $ws.Range("E44") '  +0.30%  '
Set-TextValue $ws.Range("D45") '0.0637'
Set-TextValue $ws.Range("E45") '  +2.66%  '
Set-TextValue $ws.Range("D46") '23.84'
Set-TextValue $ws.Range("E46") '  +2.22%  '
Set-TextValue $ws.Range("E47") '  +2.35%  '
Set-TextValue $ws.Range("D48") '0.0269'
Set-TextValue $ws.Range("E48") '  +1.93%  '
Set-TextValue $ws.Range("D49") '0.658'
Set-TextValue $ws.Range("E49") '  +0.23%  '
Set-TextValue $ws.Range("D50") '21.01'
Set-TextValue $ws.Range("E50") '  +6.51%  '
Set-TextValue $ws.Range("D51") '0.0996'
Set-TextValue $ws.Range("E51") '  +1.11%  '
